$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper approach: for "Price" (D) cells we force text type (NumberFormat "@")
# so numeric-looking strings like "381.72" are NOT auto-converted to numbers,
# then ClearFormats() to drop back to the default (unstyled) cell format --
# matching the source data which stores these as plain inline strings.

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.529.59"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.03%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.987.36"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.59%  "

# Row 4: TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.02%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.61%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.01"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.99%  "

# Row 7: XRP
$ws.Range("E7").Value = "  +1.17%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.00%  "

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.593"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.75%  "

# Row 10: Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.70"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.09%  "

# Row 11: TRON
$ws.Range("E11").Value = "  -0.80%  "

# Row 12: Dogecoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0857"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.88%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.454.94"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.41%  "

# Row 15: Polkadot
$ws.Range("E15").Value = "  +3.01%  "

# Row 16: WrappedEther
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.996.42"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.94%  "

# Row 17: Uniswap
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "11.17"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.39%  "

# Row 18: Polygon
$ws.Range("E18").Value = "  +0.21%  "

# Row 19: WrappedBTC
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.514.82"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.07%  "

# Row 21: InternetComputer(DFINITY)
$ws.Range("E21").Value = "  +1.18%  "

# Row 23: Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.32"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.32%  "

# Row 24: BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.69%  "

# Row 25: PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.24"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.64%  "

# Row 26: Filecoin
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.82"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.65%  "

# Row 27: RenderToken
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.59%  "

# Row 28: Kaspa
$ws.Range("E28").Value = "  +3.00%  "

# Row 29: Dai
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.07%  "

# Row 30: EthereumClassic
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.10"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.84%  "

# Row 31: Hedera
$ws.Range("E31").Value = "  -0.43%  "

# Row 32: Cosmos
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.40"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.07%  "

# Row 33: InjectiveProtocol
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.88"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.41%  "

# Row 34: OKB
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "51.39"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.30%  "

# Row 35: Toncoin
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.22%  "

# Row 36: VeChain
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0443"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.39%  "

# Row 37: FirstDigitalUSD
$ws.Range("E37").Value = "  +0.05%  "

# Row 38: LidoDAOToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.28"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.10%  "

# Row 39: Celestia
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.89"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.59%  "

# Row 40: Stacks
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.59"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.80%  "

# Row 41: Stellar
$ws.Range("E41").Value = "  +0.90%  "

# Row 42: ARBITRUM
$ws.Range("E42").Value = "  +2.86%  "

# Row 43: Monero
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "127.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.61%  "

# Row 44: NEARProtocol
$ws.Range("E44").Value = "  +12.74%  "

# Row 45: EnergySwap
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.50"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.08%  "

# Row 46: WEMIXToken
$ws.Range("E46").Value = "  +0.29%  "

# Row 47: TheGraph
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.274"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.57%  "

# Row 48: ApeXProtocol
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.37"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.66%  "

# Row 49: Maker
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.029.91"
$ws.Range("D49").ClearFormats()

# Row 50: RocketPoolETH
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.282.03"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.38%  "

# Row 51: BEAM
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0336"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.88%  "
